$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 195
$ws.Range('B195').Value = 6920523
$ws.Range('F195').Value = 'SV DrochtersenAssel'
$ws.Range('G195').Value = 'Bremer SV'
$ws.Range('H195').Value = 2
$ws.Range('I195').Value = 1
$ws.Range('K195').Value = 1.333
$ws.Range('L195').Value = 5
$ws.Range('M195').Value = 6
$ws.Range('N195').Value = 1.4
$ws.Range('O195').Value = 4.5
$ws.Range('P195').Value = 6
$ws.Range('Q195').Value = -1.25
$ws.Range('R195').Value = 1.975
$ws.Range('S195').Value = 1.875
$ws.Range('T195').Value = 2.5
$ws.Range('U195').Value = 1.95
$ws.Range('V195').Value = 1.9
$ws.Range('W195').Value = 0.3999999999999999
$ws.Range('Z195').Value = -0.5
$ws.Range('AA195').Value = 0.4375
$ws.Range('AB195').Value = 0.95
$ws.Range('AC195').Value = -1

# Row 196
$ws.Range('B196').Value = 6920524
$ws.Range('F196').Value = 'SC Weiche Flensburg 08'
$ws.Range('G196').Value = 'Eimsbutteler TV'
$ws.Range('H196').Value = 1
$ws.Range('I196').Value = 0
$ws.Range('K196').Value = 1.615
$ws.Range('L196').Value = 4
$ws.Range('M196').Value = 4
$ws.Range('N196').Value = 1.5
$ws.Range('O196').Value = 4
$ws.Range('P196').Value = 5
$ws.Range('Q196').Value = -1
$ws.Range('R196').Value = 1.85
$ws.Range('S196').Value = 1.95
$ws.Range('T196').Value = 3
$ws.Range('U196').Value = 1.825
$ws.Range('V196').Value = 1.975
$ws.Range('W196').Value = 0.5
$ws.Range('Z196').Value = 0
$ws.Range('AA196').Value = -0
$ws.Range('AB196').Value = -1
$ws.Range('AC196').Value = 0.9750000000000001

# Row 201
$ws.Range('B201').Value = 6920528
$ws.Range('F201').Value = 'Eintracht Norderstedt'
$ws.Range('G201').Value = 'Hannover II'
$ws.Range('H201').Value = 1
$ws.Range('I201').Value = 3
$ws.Range('K201').Value = 3.1
$ws.Range('M201').Value = 1.833
$ws.Range('N201').Value = 4.5
$ws.Range('O201').Value = 4.5
$ws.Range('P201').Value = 1.5
$ws.Range('Q201').Value = 1.25
$ws.Range('R201').Value = 1.825
$ws.Range('S201').Value = 1.975
$ws.Range('T201').Value = 3.5
$ws.Range('U201').Value = 2
$ws.Range('V201').Value = 1.8
$ws.Range('Y201').Value = 0.5
$ws.Range('AA201').Value = 0.9750000000000001
$ws.Range('AB201').Value = 1

# Row 202
$ws.Range('B202').Value = 6920527
$ws.Range('F202').Value = 'SV Meppen'
$ws.Range('G202').Value = 'St Pauli II'
$ws.Range('H202').Value = 3
$ws.Range('I202').Value = 4
$ws.Range('J202').Value = 'A'
$ws.Range('K202').Value = 1.6
$ws.Range('M202').Value = 4.2
$ws.Range('N202').Value = 1.75
$ws.Range('P202').Value = 3.8
$ws.Range('Q202').Value = -0.75
$ws.Range('R202').Value = 1.975
$ws.Range('S202').Value = 1.825
$ws.Range('T202').Value = 3
$ws.Range('U202').Value = 1.875
$ws.Range('V202').Value = 1.925
$ws.Range('W202').Value = -1
$ws.Range('Y202').Value = 2.8
$ws.Range('Z202').Value = -1
$ws.Range('AA202').Value = 0.825
$ws.Range('AB202').Value = 0.875
$ws.Range('AC202').Value = -1

# Row 203
$ws.Range('B203').Value = 6920529
$ws.Range('F203').Value = 'TSV Havelse'
$ws.Range('G203').Value = 'TuS BlauWeiss Lohne'
$ws.Range('I203').Value = 0
$ws.Range('J203').Value = 'H'
$ws.Range('K203').Value = 1.8
$ws.Range('M203').Value = 3.25
$ws.Range('N203').Value = 1.85
$ws.Range('O203').Value = 3.6
$ws.Range('P203').Value = 3.3
$ws.Range('Q203').Value = -0.5
$ws.Range('R203').Value = 1.925
$ws.Range('S203').Value = 1.925
$ws.Range('T203').Value = 2.75
$ws.Range('U203').Value = 1.975
$ws.Range('V203').Value = 1.875
$ws.Range('W203').Value = 0.8500000000000001
$ws.Range('Y203').Value = -1
$ws.Range('Z203').Value = 0.925
$ws.Range('AA203').Value = -1
$ws.Range('AB203').Value = -1
$ws.Range('AC203').Value = 0.875

# Row 221
$ws.Range('B221').Value = 7764326
$ws.Range('F221').Value = 'SC Weiche Flensburg 08'
$ws.Range('G221').Value = 'SV DrochtersenAssel'
$ws.Range('H221').Value = 2
$ws.Range('J221').Value = 'D'
$ws.Range('K221').Value = 2.875
$ws.Range('L221').Value = 3.75
$ws.Range('M221').Value = 2
$ws.Range('N221').Value = 2.7
$ws.Range('O221').Value = 3.4
$ws.Range('P221').Value = 2.3
$ws.Range('Q221').Value = 0.25
$ws.Range('R221').Value = 1.775
$ws.Range('S221').Value = 2.025
$ws.Range('T221').Value = 2.5
$ws.Range('U221').Value = 1.975
$ws.Range('V221').Value = 1.825
$ws.Range('X221').Value = 2.4
$ws.Range('Y221').Value = -1
$ws.Range('Z221').Value = 0.3875
$ws.Range('AA221').Value = -0.5
$ws.Range('AB221').Value = 0.9750000000000001
$ws.Range('AC221').Value = -1

# Row 222
$ws.Range('B222').Value = 7764328
$ws.Range('F222').Value = 'SC SpelleVenhaus'
$ws.Range('G222').Value = 'Eintracht Norderstedt'
$ws.Range('H222').Value = 0
$ws.Range('J222').Value = 'A'
$ws.Range('K222').Value = 3.75
$ws.Range('L222').Value = 4
$ws.Range('M222').Value = 1.666
$ws.Range('N222').Value = 5
$ws.Range('O222').Value = 4
$ws.Range('P222').Value = 1.55
$ws.Range('Q222').Value = 1
$ws.Range('R222').Value = 1.85
$ws.Range('S222').Value = 1.95
$ws.Range('T222').Value = 2.75
$ws.Range('U222').Value = 1.8
$ws.Range('V222').Value = 2
$ws.Range('X222').Value = -1
$ws.Range('Y222').Value = 0.55
$ws.Range('Z222').Value = -1
$ws.Range('AA222').Value = 0.95
$ws.Range('AB222').Value = -1
$ws.Range('AC222').Value = 1

# Row 228
$ws.Range('B228').Value = 6922680
$ws.Range('E228').Value = 45387.58333333334
$ws.Range('F228').Value = 'Eimsbutteler TV'
$ws.Range('G228').Value = 'FC Teutonia 05'
$ws.Range('K228').Value = 3.75
$ws.Range('M228').Value = 1.727
$ws.Range('N228').Value = 3.8
$ws.Range('O228').Value = 4
$ws.Range('P228').Value = 1.666
$ws.Range('Q228').Value = 0.75
$ws.Range('R228').Value = 1.925
$ws.Range('S228').Value = 1.925
$ws.Range('T228').Value = 2.75
$ws.Range('U228').Value = 1.8
$ws.Range('V228').Value = 2.05
